$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparativo proj-exe")

$ws.Range("B14").Value = 1111111
$ws.Range("B15").Value = 2222222
$ws.Range("B16").Value = 3333333
$ws.Range("B17").Value = 4444444
$ws.Range("B18").Value = 5555555
$ws.Range("B19").Value = 6666666
$ws.Range("B20").Value = 7777777
$ws.Range("B21").Value = 8888888
$ws.Range("B22").Value = 999999
$ws.Range("B23").Value = 10101010
$ws.Range("B24").Value = 1111111111
$ws.Range("B25").Value = 121211222
$ws.Range("B26").Value = 13131313
$ws.Range("B27").Value = 1414141414
$ws.Range("B28").Value = 15151515155
$ws.Range("B29").Value = 161616161

$ws.Activate()
$ws.Range("B29").Select()
